# "Generate Report for Archive" — update localization-status report:
#   1. Flip the "Ready for handoff" status to "In Translation" everywhere
#      it appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 — they all
#      share the same string, so every occurrence must move together).
#   2. Narrow the (now shorter) status columns to match their new content
#      (Overview E:F and the Status column on each language sheet).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$statusRangeOverview = $wsOverview.Range("E2:F4")
foreach ($cell in $statusRangeOverview.Cells) {
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# --- zh-cn sheet: column C (Status), rows 2-4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$statusRangeZhCn = $wsZhCn.Range("C2:C4")
foreach ($cell in $statusRangeZhCn.Cells) {
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
$wsZhCn.Range("C1").ColumnWidth = 12.5

# --- de-de sheet: column C (Status), rows 2-4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$statusRangeDeDe = $wsDeDe.Range("C2:C4")
foreach ($cell in $statusRangeDeDe.Cells) {
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
$wsDeDe.Range("C1").ColumnWidth = 12.5
